$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at F and G, shifting the old "Is Significant" column to H
$ws.Range("F1:G1").EntireColumn.Insert()

# New header labels
$ws.Range("F1").Value = "Observed"
$ws.Range("G1").Value = "Expected"
$ws.Range("H1").Value = "Is Significant"

# Observed values (column F) for each data row
$ws.Range("F2").Value = "[283 105] ; [ 76 142]"
$ws.Range("F3").Value = "[206 187] ; [ 70 154]"
$ws.Range("F4").Value = "[169 224] ; [ 51 171]"
$ws.Range("F5").Value = "[243 145] ; [ 93 129]"
$ws.Range("F6").Value = "[294 103] ; [123 100]"
$ws.Range("F7").Value = "[266 132] ; [107 118]"
$ws.Range("F8").Value = "[227 167] ; [ 85 139]"
$ws.Range("F9").Value = "[186 203] ; [ 65 156]"
$ws.Range("F10").Value = "[194 199] ; [ 72 150]"
$ws.Range("F11").Value = "[177 193] ; [ 65 146]"
$ws.Range("F12").Value = "[292 100] ; [133  92]"
$ws.Range("F13").Value = "[238 146] ; [103 119]"
$ws.Range("F14").Value = "[145 242] ; [ 51 168]"
$ws.Range("F15").Value = "[286 112] ; [131  94]"
$ws.Range("F16").Value = "[120 259] ; [ 41 176]"
$ws.Range("F17").Value = "[302  96] ; [141  82]"

# Expected values (column G) for each data row
$ws.Range("G2").Value = "[229.85478548 158.14521452] ; [129.14521452  88.85478548]"
$ws.Range("G3").Value = "[175.79902755 217.20097245] ; [100.20097245 123.79902755]"
$ws.Range("G4").Value = "[140.58536585 252.41463415] ; [ 79.41463415 142.58536585]"
$ws.Range("G5").Value = "[213.71803279 174.28196721] ; [122.28196721  99.71803279]"
$ws.Range("G6").Value = "[267.01451613 129.98548387] ; [149.98548387  73.01451613]"
$ws.Range("G7").Value = "[238.28892456 159.71107544] ; [134.71107544  90.28892456]"
$ws.Range("G8").Value = "[198.91262136 195.08737864] ; [113.08737864 110.91262136]"
$ws.Range("G9").Value = "[160.06393443 228.93606557] ; [ 90.93606557 130.06393443]"
$ws.Range("G10").Value = "[169.9804878 223.0195122] ; [ 96.0195122 125.9804878]"
$ws.Range("G11").Value = "[154.11359725 215.88640275] ; [ 87.88640275 123.11359725]"
$ws.Range("G12").Value = "[270.01620746 121.98379254] ; [154.98379254  70.01620746]"
$ws.Range("G13").Value = "[216.07920792 167.92079208] ; [124.92079208  97.07920792]"
$ws.Range("G14").Value = "[125.16831683 261.83168317] ; [ 70.83168317 148.16831683]"
$ws.Range("G15").Value = "[266.39807384 131.60192616] ; [150.60192616  74.39807384]"
$ws.Range("G16").Value = "[102.38087248 276.61912752] ; [ 58.61912752 158.38087248]"
$ws.Range("G17").Value = "[283.9194847 114.0805153] ; [159.0805153  63.9194847]"
